$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 67 - resale numbers snapshot for 2024-01-17 14:34:45
# Text columns: use a leading apostrophe so Excel stores them as literal
# text instead of auto-converting to a date/time/number, then clear the
# resulting "quote prefix" formatting so no extra style gets attached.
$ws.Range("A67").Value = "'2024-01-17"
$ws.Range("B67").Value = "'14:34:45"
$ws.Range("C67").Value = "'Wednesday"
$ws.Range("D67").Value = "'02"
$ws.Range("A67:D67").ClearFormats()

# Numeric columns
$ws.Range("E67").Value = 138925
$ws.Range("F67").Value = 139638
$ws.Range("G67").Value = 170911
$ws.Range("H67").Value = 148561
$ws.Range("I67").Value = -1
$ws.Range("J67").Value = 119128
$ws.Range("K67").Value = 222468
$ws.Range("L67").Value = 254657
$ws.Range("M67").Value = 184961
$ws.Range("N67").Value = 110346
$ws.Range("O67").Value = 41262
$ws.Range("P67").Value = 30940
$ws.Range("Q67").Value = 73414
$ws.Range("R67").Value = -1
$ws.Range("S67").Value = 42417
$ws.Range("T67").Value = -1
